# Add extra summary rows (SUMIF breakdowns) below the existing totals row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B130").Value = "ПИ_Б;Проектирование информационных систем;"
$ws.Range("C130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",C2:C125)"
$ws.Range("D130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",D2:D125)"
$ws.Range("E130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",E2:E125)"
$ws.Range("F130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",F2:F125)"
$ws.Range("G130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",G2:G125)"
$ws.Range("H130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",H2:H125)"
$ws.Range("I130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",I2:I125)"
$ws.Range("J130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",J2:J125)"
$ws.Range("K130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",K2:K125)"
$ws.Range("L130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",L2:L125)"
$ws.Range("M130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",M2:M125)"
$ws.Range("N130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",N2:N125)"
$ws.Range("O130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",O2:O125)"
$ws.Range("P130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",P2:P125)"
$ws.Range("Q130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",Q2:Q125)"
$ws.Range("R130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",R2:R125)"
$ws.Range("S130").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектирование информационных систем;"",S2:S125)"

$ws.Range("B131").Value = "ПИ_Б;Проектный практикум;"
$ws.Range("C131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",C2:C125)"
$ws.Range("D131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",D2:D125)"
$ws.Range("E131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",E2:E125)"
$ws.Range("F131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",F2:F125)"
$ws.Range("G131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",G2:G125)"
$ws.Range("H131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",H2:H125)"
$ws.Range("I131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",I2:I125)"
$ws.Range("J131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",J2:J125)"
$ws.Range("K131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",K2:K125)"
$ws.Range("L131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",L2:L125)"
$ws.Range("M131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",M2:M125)"
$ws.Range("N131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",N2:N125)"
$ws.Range("O131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",O2:O125)"
$ws.Range("P131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",P2:P125)"
$ws.Range("Q131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",Q2:Q125)"
$ws.Range("R131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",R2:R125)"
$ws.Range("S131").Formula = "=SUMIF(B2:B125,""ПИ_Б;Проектный практикум;"",S2:S125)"

$ws.Range("B132").Value = "ПИ_М;Интеграция систем;"
$ws.Range("C132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",C2:C125)"
$ws.Range("D132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",D2:D125)"
$ws.Range("E132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",E2:E125)"
$ws.Range("F132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",F2:F125)"
$ws.Range("G132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",G2:G125)"
$ws.Range("H132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",H2:H125)"
$ws.Range("I132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",I2:I125)"
$ws.Range("J132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",J2:J125)"
$ws.Range("K132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",K2:K125)"
$ws.Range("L132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",L2:L125)"
$ws.Range("M132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",M2:M125)"
$ws.Range("N132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",N2:N125)"
$ws.Range("O132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",O2:O125)"
$ws.Range("P132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",P2:P125)"
$ws.Range("Q132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",Q2:Q125)"
$ws.Range("R132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",R2:R125)"
$ws.Range("S132").Formula = "=SUMIF(B2:B125,""ПИ_М;Интеграция систем;"",S2:S125)"
